$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1210.093
$ws.Range("I40").Value = 1147.1613
$ws.Range("J40").Value = 1372.6666
$ws.Range("K40").Value = 1147.1613
$ws.Range("L40").Value = 1372.6666
$ws.Range("M40").Value = -972.1613
$ws.Range("N40").Value = -1722.6666

# Row 76
$ws.Range("H76").Value = 2060741.8
$ws.Range("I76").Value = 2181785.5
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 2181785.5
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2181470.5
$ws.Range("N76").Value = -3630

# Row 79
$ws.Range("H79").Value = 2060741.8
$ws.Range("I79").Value = 2181785.5
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 2181785.5
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2180693.5
$ws.Range("N79").Value = -5184

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 137
$ws.Range("H137").Value = 1460.75
$ws.Range("I137").Value = 1320.6818
$ws.Range("K137").Value = 3962.0454
$ws.Range("M137").Value = -1412.0454

# Row 138
$ws.Range("H138").Value = 3562.8306
$ws.Range("I138").Value = 1305.5405
$ws.Range("J138").Value = 7359.1816
$ws.Range("K138").Value = 3916.6215
$ws.Range("L138").Value = 22077.5448
$ws.Range("M138").Value = 1223.3785
$ws.Range("N138").Value = -32357.5448

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 645445.9
$ws.Range("I32").Value = 6630.6772
$ws.Range("J32").Value = 2059965.2
$ws.Range("K32").Value = 6630.6772
$ws.Range("L32").Value = 2059965.2
$ws.Range("M32").Value = -6343.6772
$ws.Range("N32").Value = -2060539.2

# Row 63
$ws.Range("H63").Value = 4115.385
$ws.Range("I63").Value = 3150
$ws.Range("J63").Value = 7333.3335
$ws.Range("K63").Value = 3150
$ws.Range("L63").Value = 7333.3335
$ws.Range("M63").Value = -2464
$ws.Range("N63").Value = -8705.333500000001

# Row 66
$ws.Range("H66").Value = 4115.385
$ws.Range("I66").Value = 3150
$ws.Range("J66").Value = 7333.3335
$ws.Range("K66").Value = 15750
$ws.Range("L66").Value = 36666.6675
$ws.Range("M66").Value = -12318
$ws.Range("N66").Value = -43530.6675

# Row 132
$ws.Range("H132").Value = 34518900
$ws.Range("I132").Value = 45455684
$ws.Range("J132").Value = 146158.86
$ws.Range("K132").Value = 136367052
$ws.Range("L132").Value = 438476.58
$ws.Range("M132").Value = -136364522
$ws.Range("N132").Value = -443536.58

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3166.5518
$ws.Range("I20").Value = 2661.7778
$ws.Range("J20").Value = 3992.5454
$ws.Range("K20").Value = 2661.7778
$ws.Range("L20").Value = 3992.5454
$ws.Range("M20").Value = -2414.7778
$ws.Range("N20").Value = -4486.5454

# Row 94
$ws.Range("H94").Value = 1177.9474
$ws.Range("I94").Value = 1427.7693
$ws.Range("J94").Value = 636.6667
$ws.Range("K94").Value = 1427.7693
$ws.Range("L94").Value = 636.6667
$ws.Range("M94").Value = -976.7692999999999
$ws.Range("N94").Value = -1538.6667

# Row 107
$ws.Range("H107").Value = 2032.1562
$ws.Range("I107").Value = 1249.0869
$ws.Range("J107").Value = 4033.3333
$ws.Range("K107").Value = 1249.0869
$ws.Range("L107").Value = 4033.3333
$ws.Range("M107").Value = 670.9131
$ws.Range("N107").Value = -7873.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3294860
$ws.Range("I31").Value = 4764127
$ws.Range("J31").Value = 146430.58
$ws.Range("K31").Value = 4764127
$ws.Range("L31").Value = 146430.58
$ws.Range("M31").Value = -4763832
$ws.Range("N31").Value = -147020.58

# Row 34
$ws.Range("H34").Value = 3294860
$ws.Range("I34").Value = 4764127
$ws.Range("J34").Value = 146430.58
$ws.Range("K34").Value = 4764127
$ws.Range("L34").Value = 146430.58
$ws.Range("M34").Value = -4763925
$ws.Range("N34").Value = -146834.58

# Row 99
$ws.Range("H99").Value = 1439.4736
$ws.Range("I99").Value = 1389.4193
$ws.Range("J99").Value = 1661.1428
$ws.Range("K99").Value = 1389.4193
$ws.Range("L99").Value = 1661.1428
$ws.Range("M99").Value = 108.5807
$ws.Range("N99").Value = -4657.1428

# Row 126
$ws.Range("H126").Value = 1439.4736
$ws.Range("I126").Value = 1389.4193
$ws.Range("J126").Value = 1661.1428
$ws.Range("K126").Value = 4168.257900000001
$ws.Range("L126").Value = 4983.428400000001
$ws.Range("M126").Value = -1698.257900000001
$ws.Range("N126").Value = -9923.428400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 494.75
$ws.Range("I97").Value = 493
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 1479
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -983
$ws.Range("N97").Value = -2492

# Row 122
$ws.Range("H122").Value = 406.875
$ws.Range("I122").Value = 408.33334
$ws.Range("J122").Value = 402.5
$ws.Range("K122").Value = 3675.00006
$ws.Range("L122").Value = 3622.5
$ws.Range("M122").Value = -1225.00006
$ws.Range("N122").Value = -8522.5

# Row 129
$ws.Range("H129").Value = 22223866
$ws.Range("I129").Value = 1162.5
$ws.Range("J129").Value = 30304848
$ws.Range("K129").Value = 3487.5
$ws.Range("L129").Value = 90914544
$ws.Range("M129").Value = 1512.5
$ws.Range("N129").Value = -90924544

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1664.3572
$ws.Range("I93").Value = 1745.8572
$ws.Range("J93").Value = 1582.8572
$ws.Range("K93").Value = 1745.8572
$ws.Range("L93").Value = 1582.8572
$ws.Range("M93").Value = -497.8571999999999
$ws.Range("N93").Value = -4078.8572

# Row 132
$ws.Range("H132").Value = 2979633.2
$ws.Range("I132").Value = 4169631.5
$ws.Range("J132").Value = 4637.25
$ws.Range("K132").Value = 12508894.5
$ws.Range("L132").Value = 13911.75
$ws.Range("M132").Value = -12506364.5
$ws.Range("N132").Value = -18971.75

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 29500
$ws.Range("J92").Value = 29500
$ws.Range("L92").Value = 29500
$ws.Range("N92").Value = -34492

# Row 122
$ws.Range("H122").Value = 1144.5883
$ws.Range("I122").Value = 1161.8462
$ws.Range("J122").Value = 1088.5
$ws.Range("K122").Value = 3485.5386
$ws.Range("L122").Value = 3265.5
$ws.Range("M122").Value = -1035.5386
$ws.Range("N122").Value = -8165.5

# Row 136
$ws.Range("H136").Value = 44293.26
$ws.Range("I136").Value = 50797.5
$ws.Range("J136").Value = 931.6667
$ws.Range("K136").Value = 152392.5
$ws.Range("L136").Value = 2795.0001
$ws.Range("M136").Value = -149842.5
$ws.Range("N136").Value = -7895.0001
